# Update "想去人数" (interest/attendance count) figures in column F across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets to the
# freshly re-scraped values (gh-pages data regeneration at 456a3b4).
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsPerform = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

$wsExhibit.Range("F2").Value = 7725
$wsExhibit.Range("F4").Value = 7894
$wsExhibit.Range("F5").Value = 43
$wsExhibit.Range("F7").Value = 35
$wsExhibit.Range("F8").Value = 6696
$wsExhibit.Range("F9").Value = 6696
$wsExhibit.Range("F10").Value = 3391
$wsExhibit.Range("F12").Value = 3727
$wsExhibit.Range("F14").Value = 51
$wsExhibit.Range("F15").Value = 43
$wsExhibit.Range("F16").Value = 67
$wsExhibit.Range("F20").Value = 45
$wsExhibit.Range("F21").Value = 319
$wsExhibit.Range("F23").Value = 332
$wsExhibit.Range("F24").Value = 3870
$wsExhibit.Range("F25").Value = 118
$wsExhibit.Range("F28").Value = 284
$wsExhibit.Range("F29").Value = 1494
$wsExhibit.Range("F32").Value = 2768
$wsExhibit.Range("F33").Value = 1867
$wsExhibit.Range("F34").Value = 34
$wsExhibit.Range("F36").Value = 62
$wsExhibit.Range("F37").Value = 27
$wsExhibit.Range("F38").Value = 3711
$wsExhibit.Range("F47").Value = 3
$wsExhibit.Range("F48").Value = 556
$wsExhibit.Range("F49").Value = 646
$wsExhibit.Range("F50").Value = 7
$wsPerform.Range("F6").Value = 414
$wsPerform.Range("F17").Value = 68
$wsAll.Range("F7").Value = 7725
$wsAll.Range("F9").Value = 7894
$wsAll.Range("F10").Value = 43
$wsAll.Range("F11").Value = 35
$wsAll.Range("F12").Value = 6696
$wsAll.Range("F13").Value = 3391
$wsAll.Range("F15").Value = 3727
$wsAll.Range("F17").Value = 51
$wsAll.Range("F18").Value = 43
$wsAll.Range("F19").Value = 67
$wsAll.Range("F24").Value = 319
$wsAll.Range("F25").Value = 332
$wsAll.Range("F26").Value = 3870
$wsAll.Range("F28").Value = 118
$wsAll.Range("F32").Value = 1494
$wsAll.Range("F35").Value = 2768
$wsAll.Range("F36").Value = 1867
$wsAll.Range("F37").Value = 34
$wsAll.Range("F39").Value = 62
$wsAll.Range("F47").Value = 68
$wsAll.Range("F49").Value = 556
$wsAll.Range("F50").Value = 646
